# Weekly CompStat update: new crime data collected
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/7/2025  Through  7/13/2025"

# --- Crime Complaints table updates (rows 15-31) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 11
$ws.Range("K15").Value = 22.222222222222
$ws.Range("L15").Value = 10
$ws.Range("M15").Value = 450
$ws.Range("N15").Value = 37.5

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 57.142857142857
$ws.Range("I16").Value = 82
$ws.Range("J16").Value = 68
$ws.Range("K16").Value = 20.588235294117
$ws.Range("L16").Value = 15.492957746478
$ws.Range("M16").Value = 115.789473684211
$ws.Range("N16").Value = -81.858407079646

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 7
$ws.Range("H17").Value = -61.111111111111
$ws.Range("I17").Value = 81
$ws.Range("J17").Value = 117
$ws.Range("K17").Value = -30.769230769230
$ws.Range("L17").Value = 20.895522388059
$ws.Range("M17").Value = 161.290322580645
$ws.Range("N17").Value = -10.989010989011

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -43.75
$ws.Range("I18").Value = 105
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = -21.641791044776
$ws.Range("L18").Value = 2.941176470588
$ws.Range("M18").Value = 7.142857142857
$ws.Range("N18").Value = -76.718403547671

# Row 19
$ws.Range("C19").Value = 20
$ws.Range("D19").Value = 29
$ws.Range("E19").Value = -31.034482758620
$ws.Range("F19").Value = 86
$ws.Range("G19").Value = 102
$ws.Range("H19").Value = -15.686274509803
$ws.Range("I19").Value = 620
$ws.Range("J19").Value = 604
$ws.Range("K19").Value = 2.649006622516
$ws.Range("L19").Value = -0.161030595813
$ws.Range("M19").Value = 9.347442680776
$ws.Range("N19").Value = -68.559837728194

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("C20").NumberFormat = "#,##0"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -62.5
$ws.Range("I20").Value = 20
$ws.Range("J20").Value = 22
$ws.Range("K20").Value = -9.090909090909
$ws.Range("L20").Value = -42.857142857142
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -95.744680851063

# Row 21
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -35.714285714285
$ws.Range("F21").Value = 118
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = -22.875816993464
$ws.Range("I21").Value = 919
$ws.Range("J21").Value = 957
$ws.Range("K21").Value = -3.970741901776
$ws.Range("L21").Value = 1.434878587196
$ws.Range("M21").Value = 22.207446808510
$ws.Range("N21").Value = -73.37775202781

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 2
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = -55.555555555555
$ws.Range("I22").Value = 55
$ws.Range("J22").Value = 49
$ws.Range("K22").Value = 12.244897959183
$ws.Range("L22").Value = 22.222222222222
$ws.Range("M22").Value = 44.736842105263

# Row 24
$ws.Range("C24").Value = 71
$ws.Range("D24").Value = 50
$ws.Range("E24").Value = 42
$ws.Range("F24").Value = 254
$ws.Range("G24").Value = 251
$ws.Range("H24").Value = 1.195219123505
$ws.Range("I24").Value = 1807
$ws.Range("J24").Value = 2045
$ws.Range("K24").Value = -11.638141809291
$ws.Range("L24").Value = -10.853478046374
$ws.Range("M24").Value = 96.199782844734

# Row 25
$ws.Range("C25").Value = 62
$ws.Range("D25").Value = 53
$ws.Range("E25").Value = 16.981132075471
$ws.Range("F25").Value = 230
$ws.Range("G25").Value = 250
$ws.Range("H25").Value = -8
$ws.Range("I25").Value = 1759
$ws.Range("J25").Value = 2012
$ws.Range("K25").Value = -12.574552683896
$ws.Range("L25").Value = -14.487117160914

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 66.666666666666
$ws.Range("F26").Value = 35
$ws.Range("H26").Value = -2.777777777777
$ws.Range("I26").Value = 243
$ws.Range("J26").Value = 242
$ws.Range("K26").Value = 0.413223140495
$ws.Range("L26").Value = 24.615384615384
$ws.Range("M26").Value = 81.343283582089

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 40
$ws.Range("L27").Value = 16.666666666666

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 8
$ws.Range("H28").Value = -38.461538461538
$ws.Range("I28").Value = 55
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -8.333333333333
$ws.Range("L28").Value = 19.565217391304

# Row 31
$ws.Range("L31").Value = 0
